$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 10 ("dataFolder" row) to make room for the
# new "plotsFile" property, shifting existing rows 10-14 down to 11-15.
$ws.Rows.Item(10).Insert()

$ws.Cells.Item(10, 1).Value = "plotsFile"
$ws.Cells.Item(10, 2).Value = "Plots.xlsx"
$ws.Cells.Item(10, 3).Value = 'Name of the excel file with plot definitions. Must be located in the "paramsFolder"'

# Update the selection to match the saved worksheet view.
$ws.Range("C10").Select()
